$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6104084714485225
$ws.Range("C2").Value = 0.1929830536217594
$ws.Range("D2").Value = 0.008214379226505031
$ws.Range("F2").Value = 0.5824442739066455
$ws.Range("G2").Value = 0.4356998002769217
$ws.Range("H2").Value = 0.5321681332819139
$ws.Range("I2").Value = 0.3877415763175662
$ws.Range("M2").Value = 0.893861172343648
$ws.Range("B3").Value = 0.5334565221141077
$ws.Range("C3").Value = 0.1702290875777237
$ws.Range("D3").Value = 0.007651219918592744
$ws.Range("F3").Value = 0.5672997128081363
$ws.Range("G3").Value = 0.4205357837508217
$ws.Range("H3").Value = 0.5302709733468731
$ws.Range("I3").Value = 0.3884873060433556
$ws.Range("M3").Value = 0.7953148875374296
$ws.Range("B4").Value = 0.4860848318545834
$ws.Range("C4").Value = 0.1561800408025817
$ws.Range("D4").Value = 0.007303532028423376
$ws.Range("F4").Value = 0.558530012890607
$ws.Range("G4").Value = 0.4116845254989272
$ws.Range("H4").Value = 0.5295008415172475
$ws.Range("I4").Value = 0.3893129358044156
$ws.Range("M4").Value = 0.7353469505384709
$ws.Range("B5").Value = 0.4667508605519686
$ws.Range("C5").Value = 0.1504355844599274
$ws.Range("D5").Value = 0.007161388479673292
$ws.Range("F5").Value = 0.5550884889733325
$ws.Range("G5").Value = 0.4081920592744268
$ws.Range("H5").Value = 0.5292858901985653
$ws.Range("I5").Value = 0.3897413065160222
$ws.Range("M5").Value = 0.7110387731222829
$ws.Range("B6").Value = 0.4635387180390751
$ws.Range("C6").Value = 0.1494805606319289
$ws.Range("D6").Value = 0.007137758629617963
$ws.Range("F6").Value = 0.554524986212698
$ws.Range("G6").Value = 0.40761902222782
$ws.Range("H6").Value = 0.5292561590788551
$ws.Range("I6").Value = 0.389817972692299
$ws.Range("M6").Value = 0.7070100368049452
$ws.Range("B7").Value = 0.485824205397563
$ws.Range("C7").Value = 0.1561026470763238
$ws.Range("D7").Value = 0.007301616856231874
$ws.Range("F7").Value = 0.5584830651447561
$ws.Range("G7").Value = 0.4116369626886041
$ws.Range("H7").Value = 0.5294975427239308
$ws.Range("I7").Value = 0.3893183415102115
$ws.Range("M7").Value = 0.7350186071727052
$ws.Range("B8").Value = 0.5839016822422138
$ws.Range("C8").Value = 0.1851538293472004
$ws.Range("D8").Value = 0.008020609943624635
$ws.Range("F8").Value = 0.5771120908417089
$ws.Range("G8").Value = 0.4303752362104518
$ws.Range("H8").Value = 0.5314318378388236
$ws.Range("I8").Value = 0.3879220434419537
$ws.Range("M8").Value = 0.8597663754030833
$ws.Range("B9").Value = 0.7752098721710468
$ws.Range("C9").Value = 0.2414952437605677
$ws.Range("D9").Value = 0.009414575706514938
$ws.Range("F9").Value = 0.6178829777707762
$ws.Range("G9").Value = 0.4708168023301482
$ws.Range("H9").Value = 0.5383749227036247
$ws.Range("I9").Value = 0.3881272556920408
$ws.Range("M9").Value = 1.108990485562089
$ws.Range("B10").Value = 0.9150931093384997
$ws.Range("C10").Value = 0.2824986988164255
$ws.Range("D10").Value = 0.01042794970818051
$ws.Range("F10").Value = 0.6504825209567002
$ws.Range("G10").Value = 0.5028547318597276
$ws.Range("H10").Value = 0.5454225575349199
$ws.Range("I10").Value = 0.3901082456735452
$ws.Range("M10").Value = 1.295350485141512
$ws.Range("B11").Value = 0.9785747755825582
$ws.Range("C11").Value = 0.3010660497691333
$ws.Range("D11").Value = 0.01088641507052301
$ws.Range("F11").Value = 0.6659006106100662
$ws.Range("G11").Value = 0.5179504071109733
$ws.Range("H11").Value = 0.5490571630469816
$ws.Range("I11").Value = 0.3914146614363929
$ws.Range("M11").Value = 1.380939431872932
$ws.Range("B12").Value = 1.002590726239532
$ws.Range("C12").Value = 0.3080845511533141
$ws.Range("D12").Value = 0.01105964243858182
$ws.Range("F12").Value = 0.6718246552298268
$ws.Range("G12").Value = 0.5237429473663724
$ws.Range("H12").Value = 0.5504955713287529
$ws.Range("I12").Value = 0.3919682672765674
$ws.Range("M12").Value = 1.413474866104522
$ws.Range("B13").Value = 0.9974195135242212
$ws.Range("C13").Value = 0.3065735539886703
$ws.Range("D13").Value = 0.01102235220669456
$ws.Range("F13").Value = 0.6705449864332991
$ws.Range("G13").Value = 0.522492015653711
$ws.Range("H13").Value = 0.5501830172665336
$ws.Range("I13").Value = 0.3918464091016887
$ws.Range("M13").Value = 1.406462101960187
$ws.Range("B14").Value = 0.9805510544196068
$ws.Range("C14").Value = 0.3016437192561341
$ws.Range("D14").Value = 0.01090067439325537
$ws.Range("F14").Value = 0.6663862654122426
$ws.Range("G14").Value = 0.5184254304049887
$ws.Range("H14").Value = 0.5491742554743695
$ws.Range("I14").Value = 0.3914590230277497
$ws.Range("M14").Value = 1.383613595440522
$ws.Range("B15").Value = 0.9702155782770205
$ws.Range("C15").Value = 0.2986224102535004
$ws.Range("D15").Value = 0.01082609266985912
$ws.Range("F15").Value = 0.6638500975968071
$ws.Range("G15").Value = 0.5159444788566958
$ws.Range("H15").Value = 0.5485644548630688
$ws.Range("I15").Value = 0.3912294261214129
$ws.Range("M15").Value = 1.369634709883996
$ws.Range("B16").Value = 0.910941261917344
$ws.Range("C16").Value = 0.2812835376425085
$ws.Range("D16").Value = 0.01039793536487821
$ws.Range("F16").Value = 0.6494868326388996
$ws.Range("G16").Value = 0.5018787853107227
$ws.Range("H16").Value = 0.54519368873774
$ws.Range("I16").Value = 0.390031075960465
$ws.Range("M16").Value = 1.289774070735277
$ws.Range("B17").Value = 0.8745385260229455
$ws.Range("C17").Value = 0.2706246455228154
$ws.Range("D17").Value = 0.01013461357250023
$ws.Range("F17").Value = 0.6408268028638986
$ws.Range("G17").Value = 0.4933843132880895
$ws.Range("H17").Value = 0.5432359127933637
$ws.Range("I17").Value = 0.3894001619619658
$ws.Range("M17").Value = 1.240995878345572
$ws.Range("B18").Value = 0.8535864145482606
$ws.Range("C18").Value = 0.264485918019659
$ws.Range("D18").Value = 0.009982921519334553
$ws.Range("F18").Value = 0.6359011088653972
$ws.Range("G18").Value = 0.4885475309228298
$ws.Range("H18").Value = 0.5421501679670371
$ws.Range("I18").Value = 0.3890753840829078
$ws.Range("M18").Value = 1.213015977732894
$ws.Range("B19").Value = 0.8464899880306689
$ws.Range("C19").Value = 0.2624060824649916
$ws.Range("D19").Value = 0.009931521243824903
$ws.Range("F19").Value = 0.6342428272959069
$ws.Range("G19").Value = 0.4869182626271993
$ws.Range("H19").Value = 0.5417894647020347
$ws.Range("I19").Value = 0.388971945982135
$ws.Range("M19").Value = 1.203555285280785
$ws.Range("B20").Value = 0.8784151408752905
$ws.Range("C20").Value = 0.2717601343211413
$ws.Range("D20").Value = 0.0101626692356156
$ws.Range("F20").Value = 0.6417429449968068
$ws.Range("G20").Value = 0.4942834833335468
$ws.Range("H20").Value = 0.5434401457171987
$ws.Range("I20").Value = 0.3894633754102017
$ws.Range("M20").Value = 1.2461804755076
$ws.Range("B21").Value = 0.9855063693745478
$ws.Range("C21").Value = 0.3030920740922625
$ws.Range("D21").Value = 0.01093642468708111
$ws.Range("F21").Value = 0.6676054533170088
$ws.Range("G21").Value = 0.5196178095747825
$ws.Range("H21").Value = 0.5494688654798381
$ws.Range("I21").Value = 0.3915712043732782
$ws.Range("M21").Value = 1.390321305978503
$ws.Range("B22").Value = 1.055360965712737
$ws.Range("C22").Value = 0.3234960196451766
$ws.Range("D22").Value = 0.01143987390219081
$ws.Range("F22").Value = 0.6850071795904
$ws.Range("G22").Value = 0.5366195707941301
$ws.Range("H22").Value = 0.5537709290274506
$ws.Range("I22").Value = 0.3932923647542594
$ws.Range("M22").Value = 1.485256662568048
$ws.Range("B23").Value = 1.018091107712053
$ws.Range("C23").Value = 0.3126128481546857
$ws.Range("D23").Value = 0.01117138568170617
$ws.Range("F23").Value = 0.6756735734754216
$ws.Range("G23").Value = 0.5275043746058543
$ws.Range("H23").Value = 0.5514415758454447
$ws.Range("I23").Value = 0.3923421066194308
$ws.Range("M23").Value = 1.434518333119243
$ws.Range("B24").Value = 0.8766625970000632
$ws.Range("C24").Value = 0.2712468134175765
$ws.Range("D24").Value = 0.01014998621885965
$ws.Range("F24").Value = 0.6413285919005318
$ws.Range("G24").Value = 0.4938768228629584
$ws.Range("H24").Value = 0.543347688056059
$ws.Range("I24").Value = 0.3894346784732221
$ws.Range("M24").Value = 1.243836322458449
$ws.Range("B25").Value = 0.7235703261313233
$ws.Range("C25").Value = 0.2263214356172512
$ws.Range("D25").Value = 0.009039289877090795
$ws.Range("F25").Value = 0.606393250757165
$ws.Range("G25").Value = 0.4594731483316679
$ws.Range("H25").Value = 0.5361567054087288
$ws.Range("I25").Value = 0.3877529350223199
$ws.Range("M25").Value = 1.041030644497908
